# Add "Formation" header to column I, row 1, and select the new cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "Formation"
$ws.Range("I1").Select() | Out-Null
